$d = $word.ActiveDocument

# --- Paragraph 16: "The project is hosted on GitHub..." -> remove lastRenderedPageBreak ---
$p16 = $d.Paragraphs(16)
$frag16 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00600A76" w:rsidRDefault="00600A76" w:rsidP="00600A76"><w:r><w:t>The project is hosted on GitHub, with a small description of how to use it.</w:t></w:r></w:p>'
$p16.Range.InsertXML($frag16)

# --- Paragraph 12: "Low-pass filter..." -> add lastRenderedPageBreak ---
$p12 = $d.Paragraphs(12)
$frag12 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00600A76" w:rsidRDefault="00600A76" w:rsidP="00600A76"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Low-pass filter your samples to turn the stair-step waves created by step 1 back into smooth curves</w:t></w:r></w:p>'
$p12.Range.InsertXML($frag12)

# --- Paragraph 4: "Accessing external DLL functions..." -> merge adjacent runs ---
$p4 = $d.Paragraphs(4)
$frag4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="006B5BCD" w:rsidRDefault="006B5BCD"><w:r><w:t xml:space="preserve">Accessing external </w:t></w:r><w:r w:rsidR="00227C83"><w:t xml:space="preserve">DLL </w:t></w:r><w:r><w:t>functions from within C# code by importing them and translating the arguments.</w:t></w:r><w:r w:rsidR="00227C83"><w:t xml:space="preserve"> Learning how to interact with native code through C# imports was a pain in the butt, </w:t></w:r><w:r w:rsidR="00227C83" w:rsidRPr="00227C83"><w:t>and</w:t></w:r><w:r w:rsidR="00227C83"><w:t xml:space="preserve"> also a positive learning experience.  </w:t></w:r></w:p>'
$p4.Range.InsertXML($frag4)

# --- Paragraph 1: split intro paragraph into three paragraphs ---
$p1 = $d.Paragraphs(1)
$fragC = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00283655" w:rsidRDefault="006B5BCD"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:t xml:space="preserve">By the end of this project I had my first </w:t></w:r><w:r><w:t>programming experience</w:t></w:r><w:r><w:t xml:space="preserve"> in:</w:t></w:r></w:p>'
$p1.Range.InsertXML($fragC)

# Insert the two new paragraphs before the (now updated) first paragraph
$p1again = $d.Paragraphs(1)
$p1again.Range.InsertParagraphBefore()
$newFirst = $d.Paragraphs(1)
$fragAB = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">This was a term-long school assignment, running during fall of 2015.  It was an opportunity to explore a wide range of concepts, many of which were completely new to me.  </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The complete program allowed a user to open, create, edit and save audio files.  </w:t></w:r><w:r><w:t>W</w:t></w:r><w:r><w:t xml:space="preserve">ithin the application you can view the </w:t></w:r><w:r><w:t>waveform and frequencies</w:t></w:r><w:r><w:t xml:space="preserve">, and select time ranges, and then cut, copy and paste samples.  </w:t></w:r><w:r><w:t xml:space="preserve">With some samples selected you can amplify them, reverse them, and even change the pitch for some fun effects.  Overall you can change sampling rate and bit rate too. </w:t></w:r></w:p>'
$newFirst.Range.InsertXML($fragAB)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
